$wb = $excel.ActiveWorkbook

# Update last-updated timestamp on Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 09:52 AM"

# Update Stock List sheet with refreshed data (2 new rows inserted at top,
# bottom 2 rows dropped, existing rows shifted down by 2)
$ws = $wb.Worksheets.Item("Stock List")

$ws.Cells.Item(2, 2).Value = "MIDWESTLTD"
$ws.Cells.Item(2, 3).Value = "MIDWESTLTD"
$ws.Cells.Item(2, 4).Value = 1117.2
$ws.Cells.Item(2, 5).Value = -1.4032
$ws.Cells.Item(2, 8).Value = 4039.8864
$ws.Cells.Item(3, 2).Value = "CAPTRU-RE1"
$ws.Cells.Item(3, 3).Value = "CAPTRU-RE1"
$ws.Cells.Item(3, 4).Value = 5.67
$ws.Cells.Item(3, 5).Value = -11.9565
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(4, 2).Value = "NIFTYCASE"
$ws.Cells.Item(4, 3).Value = "NIFTYCASE"
$ws.Cells.Item(4, 4).Value = 10.19
$ws.Cells.Item(4, 5).Value = -0.5854
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(5, 2).Value = "MOMENTUM30"
$ws.Cells.Item(5, 3).Value = "MOMENTUM30"
$ws.Cells.Item(5, 4).Value = 31.54
$ws.Cells.Item(5, 5).Value = -0.6614
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(6, 2).Value = "CANHLIFE"
$ws.Cells.Item(6, 3).Value = "CANHLIFE"
$ws.Cells.Item(6, 4).Value = 118.46
$ws.Cells.Item(6, 5).Value = 0.6286
$ws.Cells.Item(6, 8).Value = 11253.7
$ws.Cells.Item(7, 2).Value = "FLEXIADD"
$ws.Cells.Item(7, 3).Value = "FLEXIADD"
$ws.Cells.Item(7, 4).Value = 10.64
$ws.Cells.Item(7, 5).Value = -1.0233
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(8, 2).Value = "MOENERGY"
$ws.Cells.Item(8, 3).Value = "MOENERGY"
$ws.Cells.Item(8, 4).Value = 36.3
$ws.Cells.Item(8, 5).Value = -0.6568000000000001
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(9, 2).Value = "MONIFTY100"
$ws.Cells.Item(9, 3).Value = "MONIFTY100"
$ws.Cells.Item(9, 4).Value = 26.49
$ws.Cells.Item(9, 5).Value = 0.3409
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(10, 2).Value = "RUBICON"
$ws.Cells.Item(10, 3).Value = "RUBICON"
$ws.Cells.Item(10, 4).Value = 652.65
$ws.Cells.Item(10, 5).Value = -0.1453
$ws.Cells.Item(10, 8).Value = 10752.4289
$ws.Cells.Item(11, 2).Value = "CRAMC"
$ws.Cells.Item(11, 3).Value = "CRAMC"
$ws.Cells.Item(11, 4).Value = 317.2
$ws.Cells.Item(11, 5).Value = 2.3226
$ws.Cells.Item(11, 8).Value = 6325.5208
$ws.Cells.Item(12, 2).Value = "LGEINDIA"
$ws.Cells.Item(12, 3).Value = "LGEINDIA"
$ws.Cells.Item(12, 4).Value = 1633.4
$ws.Cells.Item(12, 5).Value = -0.946
$ws.Cells.Item(12, 8).Value = 110870.6825
$ws.Cells.Item(13, 2).Value = "TATACAP"
$ws.Cells.Item(13, 3).Value = "TATACAP"
$ws.Cells.Item(13, 4).Value = 329.3
$ws.Cells.Item(13, 5).Value = 0.1521
$ws.Cells.Item(13, 8).Value = 139783.5374
$ws.Cells.Item(14, 2).Value = "ELIQUID"
$ws.Cells.Item(14, 3).Value = "ELIQUID"
$ws.Cells.Item(14, 4).Value = 1004.85
$ws.Cells.Item(14, 5).Value = 0.0408
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(15, 2).Value = "WEWORK"
$ws.Cells.Item(15, 3).Value = "WEWORK"
$ws.Cells.Item(15, 4).Value = 632.15
$ws.Cells.Item(15, 5).Value = -2.4008
$ws.Cells.Item(15, 8).Value = 8472.2803
$ws.Cells.Item(16, 2).Value = "GROWWRLTY"
$ws.Cells.Item(16, 3).Value = "GROWWRLTY"
$ws.Cells.Item(16, 4).Value = 10.8
$ws.Cells.Item(16, 5).Value = -0.4608
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(17, 2).Value = "ADVANCE"
$ws.Cells.Item(17, 3).Value = "ADVANCE"
$ws.Cells.Item(17, 4).Value = 130.05
$ws.Cells.Item(17, 5).Value = -5.2666
$ws.Cells.Item(17, 8).Value = 836.0358
$ws.Cells.Item(18, 2).Value = "OMFREIGHT"
$ws.Cells.Item(18, 3).Value = "OMFREIGHT"
$ws.Cells.Item(18, 4).Value = 88.90000000000001
$ws.Cells.Item(18, 5).Value = -0.5926
$ws.Cells.Item(18, 8).Value = 299.3747
$ws.Cells.Item(19, 2).Value = "GLOTTIS"
$ws.Cells.Item(19, 3).Value = "GLOTTIS"
$ws.Cells.Item(19, 4).Value = 72.73999999999999
$ws.Cells.Item(19, 5).Value = -0.8587
$ws.Cells.Item(19, 8).Value = 672.1394
$ws.Cells.Item(20, 2).Value = "FABTECH"
$ws.Cells.Item(20, 3).Value = "FABTECH"
$ws.Cells.Item(20, 4).Value = 237.72
$ws.Cells.Item(20, 5).Value = 0.4734
$ws.Cells.Item(20, 8).Value = 1056.6843
$ws.Cells.Item(21, 2).Value = "PACEDIGITK"
$ws.Cells.Item(21, 3).Value = "PACEDIGITK"
$ws.Cells.Item(21, 4).Value = 218.85
$ws.Cells.Item(21, 5).Value = 0.1327
$ws.Cells.Item(21, 8).Value = 4723.9063
$ws.Cells.Item(22, 2).Value = "JAINREC"
$ws.Cells.Item(22, 3).Value = "JAINREC"
$ws.Cells.Item(22, 4).Value = 377.25
$ws.Cells.Item(22, 5).Value = 1.2208
$ws.Cells.Item(22, 8).Value = 13018.3623
$ws.Cells.Item(23, 2).Value = "EPACKPEB"
$ws.Cells.Item(23, 3).Value = "EPACKPEB"
$ws.Cells.Item(23, 4).Value = 301.45
$ws.Cells.Item(23, 5).Value = 1.979
$ws.Cells.Item(23, 8).Value = 3028.1254
$ws.Cells.Item(24, 2).Value = "BMWVENTLTD"
$ws.Cells.Item(24, 3).Value = "BMWVENTLTD"
$ws.Cells.Item(24, 4).Value = 69.25
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 8).Value = 600.5014
$ws.Cells.Item(25, 2).Value = "STYL"
$ws.Cells.Item(25, 3).Value = "STYL"
$ws.Cells.Item(25, 4).Value = 372.4
$ws.Cells.Item(25, 5).Value = -0.8388
$ws.Cells.Item(25, 8).Value = 6025.649
$ws.Cells.Item(26, 2).Value = "JARO"
$ws.Cells.Item(26, 3).Value = "JARO"
$ws.Cells.Item(26, 4).Value = 621.5
$ws.Cells.Item(26, 5).Value = -1.4821
$ws.Cells.Item(26, 8).Value = 1377.0134
$ws.Cells.Item(27, 2).Value = "SOLARWORLD"
$ws.Cells.Item(27, 3).Value = "SOLARWORLD"
$ws.Cells.Item(27, 4).Value = 309.1
$ws.Cells.Item(27, 5).Value = -0.6269
$ws.Cells.Item(27, 8).Value = 2679.0517
$ws.Cells.Item(28, 2).Value = "ARSSBL"
$ws.Cells.Item(28, 3).Value = "ARSSBL"
$ws.Cells.Item(28, 4).Value = 537.3
$ws.Cells.Item(28, 5).Value = 4.7266
$ws.Cells.Item(28, 8).Value = 3370.2277
$ws.Cells.Item(29, 2).Value = "GANESHCP"
$ws.Cells.Item(29, 3).Value = "GANESHCP"
$ws.Cells.Item(29, 4).Value = 274.4
$ws.Cells.Item(29, 5).Value = -2.7984
$ws.Cells.Item(29, 8).Value = 1108.9312
$ws.Cells.Item(30, 2).Value = "ATLANTAELE"
$ws.Cells.Item(30, 3).Value = "ATLANTAELE"
$ws.Cells.Item(30, 4).Value = 1003.05
$ws.Cells.Item(30, 5).Value = -1.7436
$ws.Cells.Item(30, 8).Value = 7713.116
$ws.Cells.Item(31, 2).Value = "GKENERGY"
$ws.Cells.Item(31, 3).Value = "GKENERGY"
$ws.Cells.Item(31, 4).Value = 213.85
$ws.Cells.Item(31, 5).Value = -0.7933
$ws.Cells.Item(31, 8).Value = 4337.2472
$ws.Cells.Item(32, 2).Value = "SAATVIKGL"
$ws.Cells.Item(32, 3).Value = "SAATVIKGL"
$ws.Cells.Item(32, 4).Value = 528.2
$ws.Cells.Item(32, 5).Value = -1.3079
$ws.Cells.Item(32, 8).Value = 6713.6863
$ws.Cells.Item(33, 2).Value = "IVALUE"
$ws.Cells.Item(33, 3).Value = "IVALUE"
$ws.Cells.Item(33, 4).Value = 281.45
$ws.Cells.Item(33, 5).Value = -0.3364
$ws.Cells.Item(33, 8).Value = 1506.8799
$ws.Cells.Item(34, 2).Value = "VMSTMT"
$ws.Cells.Item(34, 3).Value = "VMSTMT"
$ws.Cells.Item(34, 4).Value = 70.03
$ws.Cells.Item(34, 5).Value = -0.9056
$ws.Cells.Item(34, 8).Value = 347.5674
$ws.Cells.Item(35, 2).Value = "EUROPRATIK"
$ws.Cells.Item(35, 3).Value = "EUROPRATIK"
$ws.Cells.Item(35, 4).Value = 321.75
$ws.Cells.Item(35, 5).Value = 0.8147
$ws.Cells.Item(35, 8).Value = 3288.285
$ws.Cells.Item(36, 2).Value = "SHRINGARMS"
$ws.Cells.Item(36, 3).Value = "SHRINGARMS"
$ws.Cells.Item(36, 4).Value = 229.31
$ws.Cells.Item(36, 5).Value = -1.2616
$ws.Cells.Item(36, 8).Value = 2211.284
$ws.Cells.Item(37, 2).Value = "DEVX"
$ws.Cells.Item(37, 3).Value = "DEVX"
$ws.Cells.Item(37, 4).Value = 44.53
$ws.Cells.Item(37, 5).Value = -0.3803
$ws.Cells.Item(37, 8).Value = 401.605
$ws.Cells.Item(38, 2).Value = "URBANCO"
$ws.Cells.Item(38, 3).Value = "URBANCO"
$ws.Cells.Item(38, 4).Value = 148.9
$ws.Cells.Item(38, 5).Value = -2.0459
$ws.Cells.Item(38, 8).Value = 21380.5798
$ws.Cells.Item(39, 2).Value = "SML100CASE"
$ws.Cells.Item(39, 3).Value = "SML100CASE"
$ws.Cells.Item(39, 4).Value = 10.36
$ws.Cells.Item(39, 5).Value = -0.7663
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(40, 2).Value = "AONEGOLD"
$ws.Cells.Item(40, 3).Value = "AONEGOLD"
$ws.Cells.Item(40, 4).Value = 11.28
$ws.Cells.Item(40, 5).Value = -0.2653
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(41, 2).Value = "ELM250"
$ws.Cells.Item(41, 3).Value = "ELM250"
$ws.Cells.Item(41, 4).Value = 16.72
$ws.Cells.Item(41, 5).Value = 0.1797
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(42, 2).Value = "AMANTA"
$ws.Cells.Item(42, 3).Value = "AMANTA"
$ws.Cells.Item(42, 4).Value = 122.52
$ws.Cells.Item(42, 5).Value = 1.407
$ws.Cells.Item(42, 8).Value = 475.7372
$ws.Cells.Item(43, 2).Value = "CPEDU"
$ws.Cells.Item(43, 3).Value = "CPEDU"
$ws.Cells.Item(43, 4).Value = 315.9
$ws.Cells.Item(43, 5).Value = 1.8539
$ws.Cells.Item(43, 8).Value = 574.7148999999999
$ws.Cells.Item(44, 2).Value = "AHCL"
$ws.Cells.Item(44, 3).Value = "AHCL"
$ws.Cells.Item(44, 4).Value = 139.27
$ws.Cells.Item(44, 5).Value = 3.1706
$ws.Cells.Item(44, 8).Value = 740.2409
$ws.Cells.Item(45, 2).Value = "STLNETWORK"
$ws.Cells.Item(45, 3).Value = "STLNETWORK"
$ws.Cells.Item(45, 4).Value = 26.59
$ws.Cells.Item(45, 5).Value = -0.412
$ws.Cells.Item(45, 8).Value = 1297.3822
$ws.Cells.Item(46, 2).Value = "VIKRAN"
$ws.Cells.Item(46, 3).Value = "VIKRAN"
$ws.Cells.Item(46, 4).Value = 98.05
$ws.Cells.Item(46, 5).Value = -1.783
$ws.Cells.Item(46, 8).Value = 2528.8166
$ws.Cells.Item(47, 2).Value = "MANUFGBEES"
$ws.Cells.Item(47, 3).Value = "MANUFGBEES"
$ws.Cells.Item(47, 4).Value = 151.77
$ws.Cells.Item(47, 5).Value = -1.011
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(48, 2).Value = "MEIL"
$ws.Cells.Item(48, 3).Value = "MEIL"
$ws.Cells.Item(48, 4).Value = 461.15
$ws.Cells.Item(48, 5).Value = -0.7319
$ws.Cells.Item(48, 8).Value = 1274.1632
$ws.Cells.Item(49, 2).Value = "GROWWNXT50"
$ws.Cells.Item(49, 3).Value = "GROWWNXT50"
$ws.Cells.Item(49, 4).Value = 70.29000000000001
$ws.Cells.Item(49, 5).Value = -0.4109
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(50, 2).Value = "SHREEJISPG"
$ws.Cells.Item(50, 3).Value = "SHREEJISPG"
$ws.Cells.Item(50, 4).Value = 270.05
$ws.Cells.Item(50, 5).Value = -0.7899
$ws.Cells.Item(50, 8).Value = 4399.6074
$ws.Cells.Item(51, 2).Value = "GEMAROMA"
$ws.Cells.Item(51, 3).Value = "GEMAROMA"
$ws.Cells.Item(51, 4).Value = 219.52
$ws.Cells.Item(51, 5).Value = -0.876
$ws.Cells.Item(51, 8).Value = 1146.7097
$ws.Cells.Item(52, 2).Value = "PATELRMART"
$ws.Cells.Item(52, 3).Value = "PATELRMART"
$ws.Cells.Item(52, 4).Value = 219.31
$ws.Cells.Item(52, 5).Value = -1.0646
$ws.Cells.Item(52, 8).Value = 732.5069999999999
$ws.Cells.Item(53, 2).Value = "VIKRAMSOLR"
$ws.Cells.Item(53, 3).Value = "VIKRAMSOLR"
$ws.Cells.Item(53, 4).Value = 322
$ws.Cells.Item(53, 5).Value = -1.5892
$ws.Cells.Item(53, 8).Value = 11647.2884
$ws.Cells.Item(54, 2).Value = "LTGILTCASE"
$ws.Cells.Item(54, 3).Value = "LTGILTCASE"
$ws.Cells.Item(54, 4).Value = 29.67
$ws.Cells.Item(54, 5).Value = 0.2365
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(55, 2).Value = "REGAAL"
$ws.Cells.Item(55, 3).Value = "REGAAL"
$ws.Cells.Item(55, 4).Value = 89.13
$ws.Cells.Item(55, 5).Value = -0.8675
$ws.Cells.Item(55, 8).Value = 915.5742
$ws.Cells.Item(56, 2).Value = "BLUESTONE"
$ws.Cells.Item(56, 3).Value = "BLUESTONE"
$ws.Cells.Item(56, 4).Value = 711.95
$ws.Cells.Item(56, 5).Value = 0.1266
$ws.Cells.Item(56, 8).Value = 10773.2539
$ws.Cells.Item(57, 2).Value = "MOSILVER"
$ws.Cells.Item(57, 3).Value = "MOSILVER"
$ws.Cells.Item(57, 4).Value = 145.9
$ws.Cells.Item(57, 5).Value = -1.5054
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(58, 2).Value = "ALLTIME"
$ws.Cells.Item(58, 3).Value = "ALLTIME"
$ws.Cells.Item(58, 4).Value = 308.75
$ws.Cells.Item(58, 5).Value = 2.66
$ws.Cells.Item(58, 8).Value = 2022.5526
$ws.Cells.Item(59, 2).Value = "JSWCEMENT"
$ws.Cells.Item(59, 3).Value = "JSWCEMENT"
$ws.Cells.Item(59, 4).Value = 134.98
$ws.Cells.Item(59, 5).Value = -0.4793
$ws.Cells.Item(59, 8).Value = 18402.6999
$ws.Cells.Item(60, 2).Value = "SBILIQETF"
$ws.Cells.Item(60, 3).Value = "SBILIQETF"
$ws.Cells.Item(60, 4).Value = 1012.94
$ws.Cells.Item(60, 5).Value = 0.0296
$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(61, 2).Value = "HILINFRA"
$ws.Cells.Item(61, 3).Value = "HILINFRA"
$ws.Cells.Item(61, 4).Value = 77.23
$ws.Cells.Item(61, 5).Value = -0.3998
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(62, 2).Value = "GROWWPOWER"
$ws.Cells.Item(62, 3).Value = "GROWWPOWER"
$ws.Cells.Item(62, 4).Value = 10.28
$ws.Cells.Item(62, 5).Value = -0.9634
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(63, 2).Value = "LOTUSDEV"
$ws.Cells.Item(63, 3).Value = "LOTUSDEV"
$ws.Cells.Item(63, 4).Value = 177.82
$ws.Cells.Item(63, 5).Value = 0.3669
$ws.Cells.Item(63, 8).Value = 8690.485000000001
$ws.Cells.Item(64, 2).Value = "MBEL"
$ws.Cells.Item(64, 3).Value = "MBEL"
$ws.Cells.Item(64, 4).Value = 450.2
$ws.Cells.Item(64, 5).Value = -0.7714
$ws.Cells.Item(64, 8).Value = 2572.8126
$ws.Cells.Item(65, 2).Value = "LAXMIINDIA"
$ws.Cells.Item(65, 3).Value = "LAXMIINDIA"
$ws.Cells.Item(65, 4).Value = 145.62
$ws.Cells.Item(65, 5).Value = -1.1942
$ws.Cells.Item(65, 8).Value = 761.1248000000001
$ws.Cells.Item(66, 2).Value = "CPPLUS"
$ws.Cells.Item(66, 3).Value = "CPPLUS"
$ws.Cells.Item(66, 4).Value = 1322.1
$ws.Cells.Item(66, 5).Value = -0.264
$ws.Cells.Item(66, 8).Value = 15497.9053
$ws.Cells.Item(67, 2).Value = "SHANTIGOLD"
$ws.Cells.Item(67, 3).Value = "SHANTIGOLD"
$ws.Cells.Item(67, 4).Value = 241.57
$ws.Cells.Item(67, 5).Value = -1.6409
$ws.Cells.Item(67, 8).Value = 1741.6231
$ws.Cells.Item(68, 2).Value = "MOGOLD"
$ws.Cells.Item(68, 3).Value = "MOGOLD"
$ws.Cells.Item(68, 4).Value = 119.65
$ws.Cells.Item(68, 5).Value = -0.5403
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(69, 2).Value = "BRIGHOTEL"
$ws.Cells.Item(69, 3).Value = "BRIGHOTEL"
$ws.Cells.Item(69, 4).Value = 82.39
$ws.Cells.Item(69, 5).Value = -0.9855
$ws.Cells.Item(69, 8).Value = 3129.5229
$ws.Cells.Item(70, 2).Value = "INDIQUBE"
$ws.Cells.Item(70, 3).Value = "INDIQUBE"
$ws.Cells.Item(70, 4).Value = 212.64
$ws.Cells.Item(70, 5).Value = -0.7561
$ws.Cells.Item(70, 8).Value = 4465.6847
$ws.Cells.Item(71, 2).Value = "EBGNG"
$ws.Cells.Item(71, 3).Value = "EBGNG"
$ws.Cells.Item(71, 4).Value = 346.65
$ws.Cells.Item(71, 5).Value = 3.2311
$ws.Cells.Item(71, 8).Value = 3952.2092
$ws.Cells.Item(72, 2).Value = "LIQGRWBEES"
$ws.Cells.Item(72, 3).Value = "LIQGRWBEES"
$ws.Cells.Item(72, 4).Value = 1014.74
$ws.Cells.Item(72, 5).Value = 0.0246
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(73, 2).Value = "CHEMBONDCH"
$ws.Cells.Item(73, 3).Value = "CHEMBONDCH"
$ws.Cells.Item(73, 4).Value = 153.35
$ws.Cells.Item(73, 5).Value = -1.6987
$ws.Cells.Item(73, 8).Value = 412.459
$ws.Cells.Item(74, 2).Value = "GROWWNIFTY"
$ws.Cells.Item(74, 3).Value = "GROWWNIFTY"
$ws.Cells.Item(74, 4).Value = 10.29
$ws.Cells.Item(74, 5).Value = -0.3872
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(75, 2).Value = "ANTHEM"
$ws.Cells.Item(75, 3).Value = "ANTHEM"
$ws.Cells.Item(75, 4).Value = 702.25
$ws.Cells.Item(75, 5).Value = -0.1209
$ws.Cells.Item(75, 8).Value = 39439.0658
$ws.Cells.Item(76, 2).Value = "QUALITY30"
$ws.Cells.Item(76, 3).Value = "QUALITY30"
$ws.Cells.Item(76, 4).Value = 21.05
$ws.Cells.Item(76, 5).Value = -0.8945
$ws.Cells.Item(76, 8).Value = 0
